$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 290, shifting existing row 290 (and below) down to 291
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new data record
$ws.Cells.Item(290, 1).Value2  = 4
$ws.Cells.Item(290, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(290, 3).Value2  = "Los Lagos"
$ws.Cells.Item(290, 4).Value2  = 44943
$ws.Cells.Item(290, 5).Value2  = 10
$ws.Cells.Item(290, 6).Value2  = "Fruta"
$ws.Cells.Item(290, 7).Value2  = 100108
$ws.Cells.Item(290, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(290, 9).Value2  = 100108005
$ws.Cells.Item(290, 10).Value2 = "Piña"
$ws.Cells.Item(290, 11).Value2 = "Caramelo"
$ws.Cells.Item(290, 12).Value2 = "Segunda"
$ws.Cells.Item(290, 13).Value2 = 200
$ws.Cells.Item(290, 14).Value2 = 20000
$ws.Cells.Item(290, 15).Value2 = 21000
$ws.Cells.Item(290, 16).Value2 = 20500
$ws.Cells.Item(290, 17).Value2 = "`$/caja 14 unidades"
$ws.Cells.Item(290, 18).Value2 = "Ecuador"
$ws.Cells.Item(290, 19).Value2 = 1464
$ws.Cells.Item(290, 20).Value2 = 14
